# Add a new row (50) of memory-leak telemetry data to each of the 4 sheets.
$wb = $excel.ActiveWorkbook

$timeValue = [double]"45836.46416666666"

# --- Sheet 1: MID_LFT_#1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(50, 1).Value = $timeValue
$ws.Cells.Item(50, 2).Value = "0x01,0x90"
$ws.Cells.Item(50, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(50, 4).Value = "0x01,0x6C"
$ws.Cells.Item(50, 5).Value = "0x07"
$ws.Cells.Item(50, 6).Value = 400
$ws.Cells.Item(50, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(50, 8).Value = 364
$ws.Cells.Item(50, 9).Value = 7
$ws.Cells.Item(50, 1).NumberFormat = $ws.Cells.Item(49, 1).NumberFormat

# --- Sheet 2: MID_LFT_#2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(50, 1).Value = $timeValue
$ws.Cells.Item(50, 2).Value = "0x01,0x7c"
$ws.Cells.Item(50, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(50, 4).Value = "0x01,0x64"
$ws.Cells.Item(50, 5).Value = "0x19"
$ws.Cells.Item(50, 6).Value = 380
$ws.Cells.Item(50, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(50, 8).Value = 356
$ws.Cells.Item(50, 9).Value = 25
$ws.Cells.Item(50, 1).NumberFormat = $ws.Cells.Item(49, 1).NumberFormat

# --- Sheet 3: MID_PLT_#1 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(50, 1).Value = $timeValue
$ws.Cells.Item(50, 2).Value = "0x00,0x6e"
$ws.Cells.Item(50, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(50, 4).Value = "0x00,0x69"
$ws.Cells.Item(50, 5).Value = "0x15"
$ws.Cells.Item(50, 6).Value = 110
$ws.Cells.Item(50, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(50, 8).Value = 105
$ws.Cells.Item(50, 9).Value = 15
$ws.Cells.Item(50, 1).NumberFormat = $ws.Cells.Item(49, 1).NumberFormat

# --- Sheet 4: MID_PLT_#2 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(50, 1).Value = $timeValue
$ws.Cells.Item(50, 2).Value = "0x00,0x82"
$ws.Cells.Item(50, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(50, 4).Value = "0x00,0x7E"
$ws.Cells.Item(50, 5).Value = "0x9"
$ws.Cells.Item(50, 6).Value = 130
$ws.Cells.Item(50, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(50, 8).Value = 126
$ws.Cells.Item(50, 9).Value = 9
$ws.Cells.Item(50, 1).NumberFormat = $ws.Cells.Item(49, 1).NumberFormat
